$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# 1) "Requisitos de alto nivel" -> "Requisitos de riego"
Replace-Text "equisitos de alto nivel" "equisitos de riego"

# 2) "Descripcion de alto nivel del proyecto y sus limites" -> "Descripcion de alto riesto del proyecto y sus limites"
Replace-Text "escripción de alto nivel del proyecto y sus límites" "escripción de alto riesto del proyecto y sus límites"

# 3) merge "6. " + "El resumen del cronograma de hitos" into a single run
Replace-Text "6. El resumen del cronograma de hitos" "6. El resumen del cronograma de hitos"

# 4) merge "7. " + "Recursos financieros preaprobados" into a single run
Replace-Text "7. Recursos financieros preaprobados" "7. Recursos financieros preaprobados"

# 5) merge "8. " + "La lista de interesados clave" into a single run
Replace-Text "8. La lista de interesados clave" "8. La lista de interesados clave"

# 6) merge the "Identificar y documentar los " / "involucrados claves" / " y sus expectativas..." runs
Replace-Text "Identificar y documentar los involucrados claves y sus expectativas, cuyos intereses en el proyecto se puedan ver afectados positiva o negativamente" "Identificar y documentar los involucrados claves y sus expectativas, cuyos intereses en el proyecto se puedan ver afectados positiva o negativamente"

# 7) merge "...programación de " / "software le" / " sera entregada al cliente?" runs
Replace-Text "¿Qué documentación tales como manuales, instructivos, diagramas, algoritmos o programación de software le será entregada al cliente?" "¿Qué documentación tales como manuales, instructivos, diagramas, algoritmos o programación de software le será entregada al cliente?"

# 8) merge "9. " + "Los requisitos de aprobación del proyecto" into a single run
Replace-Text "9. Los requisitos de aprobación del proyecto" "9. Los requisitos de aprobación del proyecto"

# 9) merge "10. " + "Criterios de salida del proyecto" into a single run
Replace-Text "10. Criterios de salida del proyecto" "10. Criterios de salida del proyecto"

# 10) merge "11. " + "El director del proyecto asignado" into a single run
Replace-Text "11. El director del proyecto asignado" "11. El director del proyecto asignado"

# 11) merge "12. " + "El nombre y el nivel de autoridad..." into a single run
Replace-Text "12. El nombre y el nivel de autoridad del patrocinador o de quienes autorizan el acta de constitución del proyecto." "12. El nombre y el nivel de autoridad del patrocinador o de quienes autorizan el acta de constitución del proyecto."
